$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.144.24"
$ws.Range("E2").Value = "'  -1.11%  "
$ws.Range("D3").Value = "'3.581.76"
$ws.Range("E3").Value = "'  -1.76%  "
$ws.Range("E4").Value = "'  +0.12%  "
$ws.Range("D5").Value = "'577.38"
$ws.Range("E5").Value = "'  -2.97%  "
$ws.Range("D6").Value = "'186.55"
$ws.Range("E6").Value = "'  -4.55%  "
$ws.Range("D7").Value = "'3.580.54"
$ws.Range("E7").Value = "'  -1.63%  "
$ws.Range("E8").Value = "'  -4.48%  "
$ws.Range("E9").Value = "'  +0.01%  "
$ws.Range("E10").Value = "'  -0.53%  "
$ws.Range("D11").Value = "'0.651"
$ws.Range("E11").Value = "'  -3.87%  "
$ws.Range("D12").Value = "'55.12"
$ws.Range("E12").Value = "'  -5.40%  "
$ws.Range("D13").Value = "'0.0000305"
$ws.Range("E13").Value = "'  +2.74%  "
$ws.Range("D14").Value = "'9.56"
$ws.Range("E14").Value = "'  -4.22%  "
$ws.Range("D15").Value = "'4.156.35"
$ws.Range("E15").Value = "'  -1.81%  "
$ws.Range("D16").Value = "'19.69"
$ws.Range("E16").Value = "'  -3.29%  "
$ws.Range("D17").Value = "'3.575.77"
$ws.Range("E17").Value = "'  -2.00%  "
$ws.Range("D18").Value = "'70.050.19"
$ws.Range("E18").Value = "'  -1.24%  "
$ws.Range("D19").Value = "'12.64"
$ws.Range("E19").Value = "'  -1.10%  "
$ws.Range("E21").Value = "'  -3.03%  "
$ws.Range("D22").Value = "'488.20"
$ws.Range("E22").Value = "'  -0.30%  "
$ws.Range("D23").Value = "'19.03"
$ws.Range("E23").Value = "'  -2.06%  "
$ws.Range("D24").Value = "'4.92"
$ws.Range("E24").Value = "'  -6.53%  "
$ws.Range("E25").Value = "'  -2.03%  "
$ws.Range("D26").Value = "'95.06"
$ws.Range("E26").Value = "'  +3.87%  "
$ws.Range("D27").Value = "'11.79"
$ws.Range("E27").Value = "'  +3.15%  "
$ws.Range("E28").Value = "'  -6.32%  "
$ws.Range("D29").Value = "'9.34"
$ws.Range("E29").Value = "'  -2.76%  "
$ws.Range("D30").Value = "'7.81"
$ws.Range("E30").Value = "'  -1.07%  "
$ws.Range("D31").Value = "'31.62"
$ws.Range("E31").Value = "'  -3.65%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "'12.12"
$ws.Range("E32").Value = "'  -1.39%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").Value = "'66.36"
$ws.Range("E33").Value = "'  -0.08%  "
$ws.Range("E34").Value = "'  -6.43%  "
$ws.Range("D35").Value = "'573.64"
$ws.Range("E35").Value = "'  -7.04%  "
$ws.Range("D36").Value = "'3.28"
$ws.Range("E36").Value = "'  +15.46%  "
$ws.Range("D37").Value = "'0.424"
$ws.Range("E37").Value = "'  +3.05%  "
$ws.Range("D38").Value = "'38.81"
$ws.Range("E38").Value = "'  -3.74%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "'  +0.12%  "
$ws.Range("D40").Value = "'0.0₃0791"
$ws.Range("E40").Value = "'  -5.20%  "
$ws.Range("D41").Value = "'3.46"
$ws.Range("E41").Value = "'  -3.78%  "
$ws.Range("D42").Value = "'3.15"
$ws.Range("E42").Value = "'  -2.25%  "
$ws.Range("E43").Value = "'  -8.97%  "
$ws.Range("D44").Value = "'3.09"
$ws.Range("E44").Value = "'  -2.61%  "
$ws.Range("D45").Value = "'3.216.13"
$ws.Range("E45").Value = "'  -3.39%  "
$ws.Range("D46").Value = "'0.0444"
$ws.Range("E46").Value = "'  -3.26%  "
$ws.Range("D47").Value = "'3.48"
$ws.Range("E47").Value = "'  +4.66%  "
$ws.Range("B48").Value = "OceanProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
$ws.Range("D48").Value = "'1.65"
$ws.Range("E48").Value = "'  +36.88%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'9.56"
$ws.Range("E49").Value = "'  -0.76%  "
$ws.Range("E50").Value = "'  -2.36%  "
$ws.Range("E51").Value = "'  -0.02%  "
